# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.021.36"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "2.304.85"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'305.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").Value = "'97.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").Value = "'0.505"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").Value = "'35.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").Value = "'0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").Value = "'18.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.83%  "

$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("D14").Value = "'6.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.36%  "

$ws.Range("D15").Value = "2.663.38"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "2.317.48"
$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").Value = "42.883.48"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").Value = "'12.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.76%  "

$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").Value = "'6.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").Value = "'67.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").Value = "'236.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.55%  "

$ws.Range("D24").Value = "'2.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.70%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "'2.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.54%  "

$ws.Range("D27").Value = "'24.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").Value = "'166.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("D30").Value = "'9.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.76%  "

$ws.Range("D31").Value = "'33.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").Value = "'18.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.55%  "

$ws.Range("D34").Value = "'5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("D35").Value = "'4.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.09%  "

$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("D37").Value = "'0.0689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D39").Value = "'1.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").Value = "'2.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").Value = "1.998.33"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").Value = "'0.0281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'10.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.39%  "

$ws.Range("E45").Value = "  +5.55%  "

$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").Value = "'2.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("E48").Value = "  +4.50%  "

$ws.Range("D49").Value = "'53.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").Value = "2.531.35"
$ws.Range("E50").Value = "  +0.68%  "

# Row 51: coin swapped out (BitcoinSV -> Stacks) as ranking shifted
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
